$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 45752.6
$ws.Range("J93").Value = 45752.6
$ws.Range("L93").Value = 45752.6
$ws.Range("N93").Value = -50744.6
$ws.Range("H107").Value = 4913.0415
$ws.Range("I107").Value = 5733.15
$ws.Range("K107").Value = 5733.15
$ws.Range("M107").Value = -3813.15
$ws.Range("H109").Value = 34120
$ws.Range("J109").Value = 34120
$ws.Range("L109").Value = 34120
$ws.Range("N109").Value = -36894
$ws.Range("H113").Value = 2512.25
$ws.Range("I113").Value = 2566.4443
$ws.Range("J113").Value = 2442.5715
$ws.Range("K113").Value = 2566.4443
$ws.Range("L113").Value = 2442.5715
$ws.Range("M113").Value = 687.5556999999999
$ws.Range("N113").Value = -8950.5715
$ws.Range("H117").Value = 44404.5
$ws.Range("J117").Value = 44404.5
$ws.Range("L117").Value = 44404.5
$ws.Range("N117").Value = -53582.5
$ws.Range("H124").Value = 45884
$ws.Range("J124").Value = 45884
$ws.Range("L124").Value = 45884
$ws.Range("N124").Value = -55704
$ws.Range("H125").Value = 970.3077
$ws.Range("I125").Value = 1373.6666
$ws.Range("J125").Value = 624.5714
$ws.Range("K125").Value = 12362.9994
$ws.Range("L125").Value = 5621.1426
$ws.Range("M125").Value = -9902.999400000001
$ws.Range("N125").Value = -10541.1426
$ws.Range("H128").Value = 36318
$ws.Range("J128").Value = 36318
$ws.Range("L128").Value = 36318
$ws.Range("N128").Value = -46278
$ws.Range("H130").Value = 43298.4
$ws.Range("J130").Value = 43298.4
$ws.Range("L130").Value = 43298.4
$ws.Range("N130").Value = -53338.4
$ws.Range("H138").Value = 1847.5051
$ws.Range("I138").Value = 1733.2778
$ws.Range("J138").Value = 1914.9181
$ws.Range("K138").Value = 5199.8334
$ws.Range("L138").Value = 5744.754300000001
$ws.Range("M138").Value = -59.83340000000044
$ws.Range("N138").Value = -16024.7543

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2489.125
$ws.Range("I45").Value = 1775
$ws.Range("J45").Value = 3203.25
$ws.Range("K45").Value = 1775
$ws.Range("L45").Value = 3203.25
$ws.Range("M45").Value = -1398
$ws.Range("N45").Value = -3957.25
$ws.Range("H110").Value = 1056.4375
$ws.Range("I110").Value = 871.1539
$ws.Range("J110").Value = 1859.3334
$ws.Range("K110").Value = 871.1539
$ws.Range("L110").Value = 1859.3334
$ws.Range("M110").Value = 1173.8461
$ws.Range("N110").Value = -5949.3334
$ws.Range("H123").Value = 51429
$ws.Range("J123").Value = 51429
$ws.Range("L123").Value = 51429
$ws.Range("N123").Value = -61229
$ws.Range("H130").Value = 39510
$ws.Range("J130").Value = 39510
$ws.Range("L130").Value = 39510
$ws.Range("N130").Value = -49550
$ws.Range("H131").Value = 50412.25
$ws.Range("J131").Value = 50412.25
$ws.Range("L131").Value = 50412.25
$ws.Range("N131").Value = -60492.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 45401
$ws.Range("J92").Value = 45401
$ws.Range("L92").Value = 45401
$ws.Range("N92").Value = -50393
$ws.Range("H94").Value = 735.05884
$ws.Range("I94").Value = 637.25
$ws.Range("J94").Value = 822
$ws.Range("K94").Value = 637.25
$ws.Range("L94").Value = 822
$ws.Range("M94").Value = -186.25
$ws.Range("N94").Value = -1724
$ws.Range("H111").Value = 41097
$ws.Range("J111").Value = 41097
$ws.Range("L111").Value = 41097
$ws.Range("N111").Value = -49277
$ws.Range("H124").Value = 50996
$ws.Range("J124").Value = 50996
$ws.Range("L124").Value = 50996
$ws.Range("N124").Value = -60816
$ws.Range("H125").Value = 50676
$ws.Range("J125").Value = 50676
$ws.Range("L125").Value = 50676
$ws.Range("N125").Value = -60516
$ws.Range("H126").Value = 50768
$ws.Range("J126").Value = 50768
$ws.Range("L126").Value = 50768
$ws.Range("N126").Value = -60648
$ws.Range("H130").Value = 48917
$ws.Range("J130").Value = 48917
$ws.Range("L130").Value = 48917
$ws.Range("N130").Value = -58957

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49772
$ws.Range("J20").Value = 49772
$ws.Range("L20").Value = 49772
$ws.Range("N20").Value = -50244
$ws.Range("H30").Value = 49772
$ws.Range("J30").Value = 49772
$ws.Range("L30").Value = 49772
$ws.Range("N30").Value = -49954
$ws.Range("H31").Value = 4439.1865
$ws.Range("I31").Value = 1481.0769
$ws.Range("J31").Value = 6008.796
$ws.Range("K31").Value = 1481.0769
$ws.Range("L31").Value = 6008.796
$ws.Range("M31").Value = -1186.0769
$ws.Range("N31").Value = -6598.796
$ws.Range("H34").Value = 4439.1865
$ws.Range("I34").Value = 1481.0769
$ws.Range("J34").Value = 6008.796
$ws.Range("K34").Value = 1481.0769
$ws.Range("L34").Value = 6008.796
$ws.Range("M34").Value = -1279.0769
$ws.Range("N34").Value = -6412.796
$ws.Range("H107").Value = 779.3158
$ws.Range("I107").Value = 762.9231
$ws.Range("J107").Value = 814.8333
$ws.Range("K107").Value = 762.9231
$ws.Range("L107").Value = 814.8333
$ws.Range("M107").Value = 1157.0769
$ws.Range("N107").Value = -4654.8333
$ws.Range("H110").Value = 39023
$ws.Range("J110").Value = 39023
$ws.Range("L110").Value = 39023
$ws.Range("N110").Value = -47203
$ws.Range("H111").Value = 41733
$ws.Range("J111").Value = 41733
$ws.Range("L111").Value = 41733
$ws.Range("N111").Value = -49913
$ws.Range("H116").Value = 52000
$ws.Range("J116").Value = 52000
$ws.Range("L116").Value = 52000
$ws.Range("N116").Value = -61178
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H128").Value = 49772
$ws.Range("J128").Value = 49772
$ws.Range("L128").Value = 49772
$ws.Range("N128").Value = -59732
$ws.Range("H134").Value = 2280.1177
$ws.Range("I134").Value = 1362.4546
$ws.Range("J134").Value = 3962.5
$ws.Range("K134").Value = 4087.3638
$ws.Range("L134").Value = 11887.5
$ws.Range("M134").Value = -1552.3638
$ws.Range("N134").Value = -16957.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3611.432
$ws.Range("J5").Value = 1638.625
$ws.Range("L5").Value = 4915.875
$ws.Range("N5").Value = -5139.875
$ws.Range("H135").Value = 3611.432
$ws.Range("J135").Value = 1638.625
$ws.Range("L135").Value = 14747.625
$ws.Range("N135").Value = -19817.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 29789.666
$ws.Range("J52").Value = 29789.666
$ws.Range("L52").Value = 29789.666
$ws.Range("N52").Value = -30307.666
$ws.Range("H97").Value = 6779.952
$ws.Range("I97").Value = 1958
$ws.Range("K97").Value = 1958
$ws.Range("M97").Value = -1462
$ws.Range("H102").Value = 1802.1428
$ws.Range("I102").Value = 1601.3334
$ws.Range("J102").Value = 3007
$ws.Range("K102").Value = 1601.3334
$ws.Range("L102").Value = 3007
$ws.Range("M102").Value = 20.66660000000002
$ws.Range("N102").Value = -6251
$ws.Range("H113").Value = 5777.7393
$ws.Range("I113").Value = 6572.3687
$ws.Range("K113").Value = 6572.3687
$ws.Range("M113").Value = -4402.3687
$ws.Range("H114").Value = 48706
$ws.Range("J114").Value = 48706
$ws.Range("L114").Value = 48706
$ws.Range("N114").Value = -57384
$ws.Range("H116").Value = 48738
$ws.Range("J116").Value = 48738
$ws.Range("L116").Value = 48738
$ws.Range("N116").Value = -57916
$ws.Range("H119").Value = 48761
$ws.Range("J119").Value = 48761
$ws.Range("L119").Value = 48761
$ws.Range("N119").Value = -58437
$ws.Range("H130").Value = 52983.2
$ws.Range("J130").Value = 52983.2
$ws.Range("L130").Value = 52983.2
$ws.Range("N130").Value = -63023.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 43994
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 64991
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 64991
$ws.Range("M2").Value = -1888
$ws.Range("N2").Value = -65215
$ws.Range("H122").Value = 2205
$ws.Range("I122").Value = 2008.3334
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6025.0002
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3575.0002
$ws.Range("N122").Value = -12400
$ws.Range("H124").Value = 40248
$ws.Range("J124").Value = 40248
$ws.Range("L124").Value = 40248
$ws.Range("N124").Value = -50068
$ws.Range("H125").Value = 49707
$ws.Range("J125").Value = 49707
$ws.Range("L125").Value = 49707
$ws.Range("N125").Value = -59547
$ws.Range("H127").Value = 50531
$ws.Range("J127").Value = 50531
$ws.Range("L127").Value = 50531
$ws.Range("N127").Value = -60451
$ws.Range("H128").Value = 40748
$ws.Range("J128").Value = 40748
$ws.Range("L128").Value = 40748
$ws.Range("N128").Value = -50708
$ws.Range("H130").Value = 47968.832
$ws.Range("J130").Value = 47962.6
$ws.Range("L130").Value = 47962.6
$ws.Range("N130").Value = -58002.6
$ws.Range("H133").Value = 35240
$ws.Range("J133").Value = 35240
$ws.Range("L133").Value = 35240
$ws.Range("N133").Value = -40300

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 50711
$ws.Range("J128").Value = 50711
$ws.Range("L128").Value = 50711
$ws.Range("N128").Value = -60671
$ws.Range("H131").Value = 48709
$ws.Range("J131").Value = 48709
$ws.Range("L131").Value = 48709
$ws.Range("N131").Value = -58789
